$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new "team record" columns, matching the
# existing header style used by the other header cells (A1:AC1) by
# copying the format from the last existing header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill the team record values for every data row (2-50)
for ($row = 2; $row -le 50; $row++) {
    $ws.Cells.Item($row, 30).Value = 74   # AD
    $ws.Cells.Item($row, 31).Value = 88   # AE
    $ws.Cells.Item($row, 32).Value = 0    # AF
}
